# Add a "2022-Q1" sheet (positioned between "2021-Q4" and "总计") with the
# new quarter's fund-holding data, and update the "总计" (totals) sheet with
# a new summary row for 2022-Q1 (inserted above the existing 2021-Q4 row).

$wb = $excel.ActiveWorkbook
$sheetQ4    = $wb.Worksheets.Item("2021-Q4")
$sheetTotal = $wb.Worksheets.Item("总计")

# ---------------------------------------------------------------------
# 1) Insert the new "2022-Q1" worksheet right before "总计".
# ---------------------------------------------------------------------
$newSheet = $wb.Worksheets.Add($sheetTotal)
$newSheet.Name = "2022-Q1"

# NOTE: inserting a sheet "before" $sheetTotal shifts everything at/after
# that slot over by one, and this host's worksheet handles are
# position-based, so the *old* $sheetTotal variable now silently refers
# to the newly inserted sheet instead of "总计". Re-resolve it by name so
# later writes in step 5 land on the correct sheet.
$sheetTotal = $wb.Worksheets.Item("总计")

# ---------------------------------------------------------------------
# 2) Clone the formatting from "2021-Q4" (header row + styled/bordered
#    index column) so the new sheet matches the existing look.
# ---------------------------------------------------------------------
$sheetQ4.Range("B1:H1").Copy()
$newSheet.Range("B1:H1").PasteSpecial(-4122)

$sheetQ4.Range("A2:H5").Copy()
$newSheet.Range("A2:H6").PasteSpecial(-4122)

$newSheet.Range("A1").NumberFormat = "General"

# ---------------------------------------------------------------------
# 3) Header row.
# ---------------------------------------------------------------------
$newSheet.Range("B1").Value = "基金代码"
$newSheet.Range("C1").Value = "基金名称"
$newSheet.Range("D1").Value = "基金规模"
$newSheet.Range("E1").Value = "股票总仓位"
$newSheet.Range("F1").Value = "仓位占比"
$newSheet.Range("G1").Value = "持有市值(亿元)"
$newSheet.Range("H1").Value = "仓位排名"

# ---------------------------------------------------------------------
# 4) Data rows. Columns B,D,E,F,G hold numeric-looking text (fund codes
#    with leading zeros, decimal figures formatted as plain strings in
#    the source data) so they must be forced to Text before assignment;
#    otherwise Excel auto-coerces them to numbers. Column C is already
#    non-numeric text. Columns A (index) and H (rank) are real numbers.
#    (Each column is formatted with its own single-area Range call —
#    a combined multi-area union string is unreliable here.)
# ---------------------------------------------------------------------
$newSheet.Range("B2:B6").NumberFormat = "@"
$newSheet.Range("D2:D6").NumberFormat = "@"
$newSheet.Range("E2:E6").NumberFormat = "@"
$newSheet.Range("F2:F6").NumberFormat = "@"
$newSheet.Range("G2:G6").NumberFormat = "@"

$data = @(
    @(0, "013393", "信达澳银价值精选混合A",     "3.61",  "81.31", "3.06", "0.1105", 7),
    @(1, "003655", "信达澳银新财富灵活配置混合", "11.86", "25.86", "0.65", "0.0771", 7),
    @(2, "012005", "信达澳银恒盛混合A",         "1.87",  "31.90", "0.75", "0.0140", 8),
    @(3, "013394", "信达澳银价值精选混合C",     "0.37",  "81.31", "3.06", "0.0113", 7),
    @(4, "012006", "信达澳银恒盛混合C",         "0.31",  "31.90", "0.75", "0.0023", 8)
)

$r = 2
foreach ($row in $data) {
    $newSheet.Cells.Item($r, 1).Value = $row[0]
    $newSheet.Cells.Item($r, 2).Value = $row[1]
    $newSheet.Cells.Item($r, 3).Value = $row[2]
    $newSheet.Cells.Item($r, 4).Value = $row[3]
    $newSheet.Cells.Item($r, 5).Value = $row[4]
    $newSheet.Cells.Item($r, 6).Value = $row[5]
    $newSheet.Cells.Item($r, 7).Value = $row[6]
    $newSheet.Cells.Item($r, 8).Value = $row[7]
    $r = $r + 1
}

# ---------------------------------------------------------------------
# 5) Update the "总计" (totals) sheet: insert a new row above the
#    existing 2021-Q4 summary row and fill it in with the 2022-Q1 totals.
#    (Literal values are used for the shifted-down row instead of
#    reading .Value back off a cell, which this host does not support
#    reliably as an expression.)
# ---------------------------------------------------------------------
$sheetTotal.Range("A2:D2").Copy()
$sheetTotal.Range("A3:D3").PasteSpecial(-4122)

$sheetTotal.Cells.Item(3, 1).Value = 1
$sheetTotal.Cells.Item(3, 2).Value = "2021-Q4"
$sheetTotal.Cells.Item(3, 3).Value = 4
$sheetTotal.Cells.Item(3, 4).Value = 1.03

$sheetTotal.Cells.Item(2, 1).Value = 0
$sheetTotal.Cells.Item(2, 2).Value = "2022-Q1"
$sheetTotal.Cells.Item(2, 3).Value = 5
$sheetTotal.Cells.Item(2, 4).Value = 0.22
